# Append 5 new flight-departure rows (rows 63-67) to the "Main Data" sheet,
# mirroring the existing layout:
#   A: NUMBER  B: DATE  C: TIME  D: FLIGHT  E: TO  F: SHORT
#   G: AIRLINE H: MODEL I: AIRCFAT ID  J: STATUS  K: (blank)  L: DIFFERENCE  M: (blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=63; A=62; B="Monday, Jan 09"; C="2:55 PM";  D="FR8012"; E="Dublin";    F="(DUB)"; G="Ryanair "; H="B38M"; I="(EI-HES)"; J="3:42 PM"; L="0 hours, 47 minutes" },
    @{ Row=64; A=63; B="Monday, Jan 09"; C="3:35 PM";  D="FR1970"; E="Catania";   F="(CTA)"; G="Ryanair "; H="B738"; I="(SP-RKP)"; J="3:53 PM"; L="0 hours, 18 minutes" },
    @{ Row=65; A=64; B="Monday, Jan 09"; C="3:40 PM";  D="FR6945"; E="Barcelona"; F="(BCN)"; G="Ryanair "; H="B738"; I="(EI-DYC)"; J="3:49 PM"; L="0 hours, 9 minutes" },
    @{ Row=66; A=65; B="Monday, Jan 09"; C="4:10 PM";  D="FR1968"; E="Madrid";    F="(MAD)"; G="Ryanair "; H="B38M"; I="(SP-RZO)"; J="4:15 PM"; L="0 hours, 5 minutes" },
    @{ Row=67; A=66; B="Monday, Jan 09"; C="5:20 PM";  D="FR4528"; E="Oslo";      F="(TRF)"; G="Ryanair "; H="B738"; I="(SP-RSV)"; J="5:18 PM"; L="0 hours, -2 minutes" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("L$n").Value = $r.L
}
